# Update "想去人数" (interest count) figures on the two sheets that list
# every event: "展览" (sheet1) and "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6920
$wsExhibit.Range("F4").Value = 449
$wsExhibit.Range("F18").Value = 3533
$wsExhibit.Range("F22").Value = 2157
$wsExhibit.Range("F23").Value = 214

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6920
$wsAll.Range("F4").Value = 449
$wsAll.Range("F19").Value = 3533
$wsAll.Range("F23").Value = 2157
$wsAll.Range("F24").Value = 214
